$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 'products__item\ in-stock\ products__item_3-in-row"]:nth-child(2) [type="button'
$ws.Range("B4").Value = "New!iPhone 15 Silicone Case"
$ws.Range("C4").Value = "Buy"
$ws.Range("D4").Value = "\31 52184-case-685"
$ws.Range("E4").Value = "\31 52185-case-689"
$ws.Range("F4").Value = "iPhone 15 Plus"
$ws.Range("G4").Value = "Cypress"

$ws.Range("A5").Value = 'products__item\ in-stock\ products__item_3-in-row"]:nth-child(2) [type="button'
$ws.Range("B5").Value = "New!iPhone 15 Silicone Case"
$ws.Range("C5").Value = "Buy"
$ws.Range("D5").Value = "\31 52184-case-685"
$ws.Range("E5").Value = "\31 52185-case-689"
$ws.Range("F5").Value = "iPhone 15 Plus"
$ws.Range("G5").Value = "Cypress"
